$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 24 de Marzo de 2020 a las 17:46'
$ws.Range('B6').Value = 49344
$ws.Range('C6').Value = 5610
$ws.Range('D6').Value = 297
$ws.Range('E6').Value = 48431
$ws.Range('G6').Value = 63
$ws.Range('H6').Value = 616
$ws.Range('B8').Value = 31991
$ws.Range('C8').Value = 2935
$ws.Range('E8').Value = 31093
$ws.Range('G8').Value = 26
$ws.Range('H8').Value = 149
$ws.Range('B17').Value = 2768
$ws.Range('C17').Value = 143
$ws.Range('E17').Value = 2750
$ws.Range('E18').Value = 2307
$ws.Range('G18').Value = 10
$ws.Range('H18').Value = 33
$ws.Range('B20').Value = 2177
$ws.Range('C20').Value = 86
$ws.Range('E20').Value = 2041
$ws.Range('B25').Value = 1591
$ws.Range('C25').Value = 131
$ws.Range('E25').Value = 1558
$ws.Range('A34').Value = 'Polonia'
$ws.Range('B34').Value = 848
$ws.Range('C34').Value = 99
$ws.Range('D34').Value = 1
$ws.Range('E34').Value = 838
$ws.Range('F34').Value = 3
$ws.Range('G34').Value = 1
$ws.Range('H34').Value = 9
$ws.Range('A35').Value = 'Tailandia'
$ws.Range('B35').Value = 827
$ws.Range('C35').Value = 106
$ws.Range('D35').Value = 52
$ws.Range('E35').Value = 771
$ws.Range('F35').Value = 7
$ws.Range('G35').Value = 3
$ws.Range('H35').Value = 4
$ws.Range('A39').Value = 'Grecia'
$ws.Range('B39').Value = 743
$ws.Range('C39').Value = 48
$ws.Range('D39').Value = 29
$ws.Range('E39').Value = 694
$ws.Range('F39').Value = 35
$ws.Range('G39').Value = 3
$ws.Range('H39').Value = 20
$ws.Range('A40').Value = 'Crucero'
$ws.Range('B40').Value = 712
$ws.Range('D40').Value = 587
$ws.Range('E40').Value = 115
$ws.Range('F40').Value = 15
$ws.Range('H40').Value = 10
$ws.Range('B65').Value = 249
$ws.Range('C65').Value = 14
$ws.Range('D65').Value = 4
$ws.Range('E65').Value = 245
$ws.Range('A66').Value = 'Emiratos Arabes Unidos'
$ws.Range('B66').Value = 248
$ws.Range('C66').Value = 50
$ws.Range('D66').Value = 45
$ws.Range('E66').Value = 201
$ws.Range('F66').Value = 2
$ws.Range('H66').Value = 2
$ws.Range('A67').Value = 'Bulgaria'
$ws.Range('B67').Value = 218
$ws.Range('C67').Value = 17
$ws.Range('D67').Value = 3
$ws.Range('E67').Value = 212
$ws.Range('F67').Value = 8
$ws.Range('H67').Value = 3
$ws.Range('A68').Value = 'Taiwan'
$ws.Range('B68').Value = 216
$ws.Range('C68').Value = 21
$ws.Range('D68').Value = 29
$ws.Range('E68').Value = 185
$ws.Range('F68').Value = 0
$ws.Range('H68').Value = 2
$ws.Range('A69').Value = 'Eslovaquia'
$ws.Range('B69').Value = 204
$ws.Range('C69').Value = 18
$ws.Range('D69').Value = 7
$ws.Range('E69').Value = 197
$ws.Range('F69').Value = 2
$ws.Range('H69').Value = 0
$ws.Range('A70').Value = 'Lituania'
$ws.Range('B70').Value = 203
$ws.Range('C70').Value = 24
$ws.Range('D70').Value = 1
$ws.Range('E70').Value = 201
$ws.Range('F70').Value = 1
$ws.Range('H70').Value = 1
$ws.Range('E87').Value = 110
$ws.Range('G87').Value = 2
$ws.Range('H87').Value = 3
$ws.Range('A129').Value = 'Isla de Man'
$ws.Range('C129').Value = 10
$ws.Range('A130').Value = 'Polinesia Francesa'
$ws.Range('C130').Value = 5
$ws.Range('D130').Value = 0
$ws.Range('E130').Value = 23
$ws.Range('A131').Value = 'Monaco'
$ws.Range('C131').Value = 0
$ws.Range('D131').Value = 1
$ws.Range('E131').Value = 22
$ws.Range('A132').Value = 'Guayana Francesa'
$ws.Range('B132').Value = 23
$ws.Range('C132').Value = 3
$ws.Range('D132').Value = 6
$ws.Range('E132').Value = 17
$ws.Range('H132').Value = 0
$ws.Range('A133').Value = 'Jamaica'
$ws.Range('B133').Value = 21
$ws.Range('C133').Value = 2
$ws.Range('D133').Value = 2
$ws.Range('E133').Value = 18
$ws.Range('H133').Value = 1
$ws.Range('A135').Value = 'Togo'
$ws.Range('C135').Value = 2
$ws.Range('D135').Value = 1
$ws.Range('H135').Value = 0
$ws.Range('A136').Value = 'Guatemala'
$ws.Range('C136').Value = 0
$ws.Range('D136').Value = 0
$ws.Range('H136').Value = 1
$ws.Range('A138').Value = 'Barbados'
$ws.Range('A139').Value = 'Islas Virgenes de los Estados Unidos'
$ws.Range('A180').Value = 'Gambia'
$ws.Range('B180').Value = 3
$ws.Range('C180').Value = 1
$ws.Range('H180').Value = 1
$ws.Range('A184').Value = 'Butan'
$ws.Range('A185').Value = 'Birmania'
$ws.Range('A186').Value = 'Mauritania'
$ws.Range('A187').Value = 'Dominica'
$ws.Range('D187').Value = 0
$ws.Range('E187').Value = 2
$ws.Range('A188').Value = 'Nepal'
$ws.Range('D188').Value = 1
$ws.Range('H188').Value = 0
$ws.Range('A189').Value = 'Santa Sede'
$ws.Range('A190').Value = 'Islas Turcas y Caicos'
$ws.Range('A191').Value = 'Papua Nueva Guinea'
$ws.Range('A192').Value = 'Timor Oriental'
$ws.Range('A193').Value = 'Eritrea'
$ws.Range('A194').Value = 'Siria'
$ws.Range('A195').Value = 'Granada'
$ws.Range('A196').Value = 'Belice'
$ws.Range('A197').Value = 'Montserrat'
$ws.Range('A199').Value = 'Somalia'
